$d = $word.ActiveDocument

$d.Content.Find.Execute("366×4=1464", $true, $false, $false, $false, $false, $true, 1, $false, "592×4=2368", 2) | Out-Null
$d.Content.Find.Execute("867×7=6069", $true, $false, $false, $false, $false, $true, 1, $false, "181×4=724", 2) | Out-Null
$d.Content.Find.Execute("498×2=996", $true, $false, $false, $false, $false, $true, 1, $false, "521×4=2084", 2) | Out-Null
$d.Content.Find.Execute("473×9=4257", $true, $false, $false, $false, $false, $true, 1, $false, "557×8=4456", 2) | Out-Null
$d.Content.Find.Execute("297×4=1188", $true, $false, $false, $false, $false, $true, 1, $false, "364×6=2184", 2) | Out-Null
$d.Content.Find.Execute("966×2=1932", $true, $false, $false, $false, $false, $true, 1, $false, "854×7=5978", 2) | Out-Null
$d.Content.Find.Execute("413×5=2065", $true, $false, $false, $false, $false, $true, 1, $false, "967×2=1934", 2) | Out-Null
$d.Content.Find.Execute("510×2=1020", $true, $false, $false, $false, $false, $true, 1, $false, "628×7=4396", 2) | Out-Null
$d.Content.Find.Execute("578×7=4046", $true, $false, $false, $false, $false, $true, 1, $false, "620×3=1860", 2) | Out-Null
$d.Content.Find.Execute("256×2=512", $true, $false, $false, $false, $false, $true, 1, $false, "130×4=520", 2) | Out-Null
$d.Content.Find.Execute("132×5=660", $true, $false, $false, $false, $false, $true, 1, $false, "502×3=1506", 2) | Out-Null
$d.Content.Find.Execute("679×3=2037", $true, $false, $false, $false, $false, $true, 1, $false, "288×8=2304", 2) | Out-Null
$d.Content.Find.Execute("863×6=5178", $true, $false, $false, $false, $false, $true, 1, $false, "756×4=3024", 2) | Out-Null
$d.Content.Find.Execute("946×6=5676", $true, $false, $false, $false, $false, $true, 1, $false, "527×3=1581", 2) | Out-Null
$d.Content.Find.Execute("931×6=5586", $true, $false, $false, $false, $false, $true, 1, $false, "965×3=2895", 2) | Out-Null
$d.Content.Find.Execute("888×5=4440", $true, $false, $false, $false, $false, $true, 1, $false, "906×3=2718", 2) | Out-Null
$d.Content.Find.Execute("414×4=1656", $true, $false, $false, $false, $false, $true, 1, $false, "130×7=910", 2) | Out-Null
$d.Content.Find.Execute("210×3=630", $true, $false, $false, $false, $false, $true, 1, $false, "962×5=4810", 2) | Out-Null
$d.Content.Find.Execute("644×5=3220", $true, $false, $false, $false, $false, $true, 1, $false, "174×5=870", 2) | Out-Null
$d.Content.Find.Execute("821×6=4926", $true, $false, $false, $false, $false, $true, 1, $false, "471×4=1884", 2) | Out-Null
$d.Content.Find.Execute("334×8=2672", $true, $false, $false, $false, $false, $true, 1, $false, "865×2=1730", 2) | Out-Null
$d.Content.Find.Execute("305×6=1830", $true, $false, $false, $false, $false, $true, 1, $false, "979×9=8811", 2) | Out-Null
$d.Content.Find.Execute("803×3=2409", $true, $false, $false, $false, $false, $true, 1, $false, "161×2=322", 2) | Out-Null
$d.Content.Find.Execute("643×2=1286", $true, $false, $false, $false, $false, $true, 1, $false, "569×3=1707", 2) | Out-Null
$d.Content.Find.Execute("578×9=5202", $true, $false, $false, $false, $false, $true, 1, $false, "473×8=3784", 2) | Out-Null
